# Add PA10 assignment row (with its deadline) right after the existing
# PA9 row, pushing the SA* rows down by one. This mirrors inserting a new
# row 12 ("PA10" / "December 11, 2024") in the Assignment/Deadline table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 12 (the first SA row), which
# shifts SA1..SA10 down to rows 13..22 and keeps the inherited formatting
# (text style) for the new row's cells.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row with the PA10 assignment details.
$ws.Range("A12").Value = "PA10"
$ws.Range("B12").Value = "December 11, 2024"

# Match the saved selection/active cell shown in the target workbook.
$ws.Range("B13").Select()
